# Motivator-main / Users.xlsx - "Se arreglan bugs relacionados con el NewLogin"
#
# Adds two new registered users to the "Usuarios" sheet (rows 3 and 4),
# and removes the leftover "Hipervínculo" (hyperlink) character style /
# direct formatting that was still sitting on D2 (the admin e-mail cell)
# even though it visually should just look like normal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Sebastian Jerez -------------------------------------------------
$ws.Range("A3").Value = "Sebastian Jerez"
$ws.Range("B3").Value = "Starjerez"
# Leading apostrophe forces this numeric-looking password to be stored as
# text (like the existing "admin1234" cell) instead of a number.
$ws.Range("C3").Value = "'272426799"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "sebastianjs99@hotmail.com"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 2

# --- Row 4: pepito perez -----------------------------------------------------
$ws.Range("A4").Value = "pepito perez"
$ws.Range("B4").Value = "pepito"
$ws.Range("C4").Value = "'12345678"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "pepito@pepito.com"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 3

# --- Clean up the stale hyperlink look on D2 --------------------------------
# D2 keeps its mailto: hyperlink (see sheet1.xml.rels) but should no longer
# render with the blue/underlined "Hipervínculo" style.
$ws.Range("D2").Style = "Normal"

# The "Hipervínculo" cell style is no longer used anywhere, so drop it from
# the workbook's style list.
$wb.Styles.Item("Hipervínculo").Delete() | Out-Null

# --- Selection -------------------------------------------------------------
$ws.Range("F2").Select() | Out-Null
